$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (players, positions, teams) for rows 2-19, column by column,
# matching the re-uploaded data from the commit.
$players = @('Scoot Henderson', 'Desmond Bane', 'Trae Young', 'Bilal Coulibaly', 'LeBron James', 'Walker Kessler', 'Christian Braun', 'Josh Okogie', 'Kawhi Leonard', 'Norman Powell', 'Jalen Brunson', 'Devin Booker', 'Myles Turner', 'Jalen Williams', 'D''Angelo Russell', 'Immanuel Quickley', 'Brandon Ingram', 'Jimmy Butler')
$positions = @('PG', 'SG,SF', 'PG', 'SG,SF', 'SF,PF', 'C', 'SG,SF', 'SG,SF', 'SG,SF,PF', 'SG,SF', 'PG', 'PG,SG', 'C', 'SG,SF,PF,C', 'PG', 'PG,SG', 'SG,SF,PF', 'SF,PF')
$teams = @('Portland Trail Blazers', 'Memphis Grizzlies', 'Atlanta Hawks', 'Washington Wizards', 'Los Angeles Lakers', 'Utah Jazz', 'Denver Nuggets', 'Charlotte Hornets', 'LA Clippers', 'New York Knicks', 'Phoenix Suns', 'Indiana Pacers', 'Oklahoma City Thunder', 'Brooklyn Nets', 'Toronto Raptors', 'Portland Trail Blazers', 'New Orleans Pelicans', 'Miami Heat')

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $players[$i]
}
for ($i = 0; $i -lt $positions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $positions[$i]
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value2 = $teams[$i]
}
